$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Paragraph "Agile has been shown to perform better..." (para 6)
#    - "Agile" -> "Agile methodologies"
#    - "did agile outperform" -> "did agile methodologies and practices
#      outperform"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Agile has been shown", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Agile methodologies has been shown", 2) | Out-Null

$d.Content.Find.Execute("did agile outperform", $true, $false, $false, $false, $false, `
    $true, 1, $false, "did agile methodologies and practices outperform", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Paragraph "While agile is most suited..." (para 7)
#    - "While agile is most suited" -> "While scrum is most suited"
#    - "found that agile (even" -> "found that overall scrum (even"
# ------------------------------------------------------------------
$d.Content.Find.Execute("While agile is most suited", $true, $false, $false, $false, $false, `
    $true, 1, $false, "While scrum is most suited", 2) | Out-Null

$d.Content.Find.Execute("found that agile (even", $true, $false, $false, $false, $false, `
    $true, 1, $false, "found that overall scrum (even", 2) | Out-Null

# ------------------------------------------------------------------
# 3) New paragraph inserted right after the "Project Planning and
#    Tracking" heading: "Product backlogs are a part of ...".
#    Insert it *before* the following body paragraph so the new
#    paragraph picks up the Normal style of that paragraph instead of
#    the Heading2 style of the preceding heading.
# ------------------------------------------------------------------
$planningBodyPara = $d.Paragraphs.Item(14)
$planningBodyPara.Range.InsertParagraphBefore()
$backlogPara = $d.Paragraphs.Item(14)
$backlogPara.Range.Text = "Product backlogs are a part of the scrum planning process, they allow teams to prioritize features of a given software and track progress based on estimated work (hours) left according to each user story (feature) through a gaant chart."

# ------------------------------------------------------------------
# 4) Rewrite the paragraph that used to start "Although sprint
#    backlogs are a practice of scrum..." (now paragraph 15) with the
#    new text about individual practices / Lagerberg findings.
# ------------------------------------------------------------------
$rewrite1 = $d.Paragraphs.Item(15)
$rewrite1.Range.Text = "From my research I was unable to find any empirical papers to support the benefits of specific individual practices. Though what I did find was for the most part scrum does show some benefits over more traditional methodologies as in her study Lagerberg found Project B (scrum project group) rated higher or equally in all areas except productivity [1]. Lagerberg did touch on this in her paper hinting that these perceived benefits could have been a result of planning amongst other things and may be lost if less time were to be spent on this, with more research needing to be done in this area [1]. I do agree with her on both points as scrum is a methodology with a strong focus on planning and re-planning especially at the beginning of sprints. "

# ------------------------------------------------------------------
# 5) Rewrite the paragraph that used to start " From my research I was
#    unable to find any empirical papers..." (now paragraph 16) with
#    the new text about other research / Serrador & Pinto planning
#    correlation.
# ------------------------------------------------------------------
$rewrite2 = $d.Paragraphs.Item(16)
$rewrite2.Range.Text = "Of course there were many limitation in this study in that they only looked at one set of project groups within a corporation or participant error amongst other things however there is other research that shows the benefits of agile practices. Serrador & Pinto found in their research that there is some correlation between the amount of time spent planning and the success of a project [2] thus the benefits of a product backlog along with the necessary planning that goes along with it cannot be discounted. "

# Relocate the _GoBack bookmark (previously between "(even" and
# "single practices)" in paragraph 7) to sit right before "thus" in
# the rewritten paragraph above.
$bmRange = $d.Content
$bmRange.Find.Execute("thus") | Out-Null
$bmRange.Collapse(1)
$bmRange.Bookmarks.Add("_GoBack") | Out-Null

# ------------------------------------------------------------------
# 6) New empty paragraph inserted right before "Managing Change".
#    Insert it *after* the preceding paragraph so it does not inherit
#    the "Managing Change" heading style.
# ------------------------------------------------------------------
$precedingPara = $d.Paragraphs.Item(16)
$precedingPara.Range.InsertParagraphAfter()

Write-Output "done"
